# Auto-generated edit script: update crypto price (D) and 1h volume change (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.528.30"
$ws.Range("E2").Value = "  -2.66%  "
$ws.Range("D3").Value = "1.815.24"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.70%  "
$ws.Range("D5").Value = "'308.69"
$ws.Range("E5").Value = "  -1.71%  "
$ws.Range("E6").Value = "  +0.55%  "
$ws.Range("D7").Value = "'0.4568"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("D8").Value = "'0.3668"
$ws.Range("D9").Value = "'0.07136"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "'0.8811"
$ws.Range("E10").Value = "  -1.08%  "
$ws.Range("D11").Value = "'0.07750"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "'19.36"
$ws.Range("E12").Value = "  -3.54%  "
$ws.Range("D13").Value = "1.844.98"
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "'5.303"
$ws.Range("E14").Value = "  -1.96%  "
$ws.Range("D15").Value = "'6.375"
$ws.Range("E15").Value = "  -2.11%  "
$ws.Range("D16").Value = "'86.79"
$ws.Range("E16").Value = "  -5.39%  "
$ws.Range("D17").Value = "'1.009"
$ws.Range("E17").Value = "  +0.71%  "
$ws.Range("D18").Value = "'0.000008594"
$ws.Range("E18").Value = "  -3.64%  "
$ws.Range("E19").Value = "  +0.61%  "
$ws.Range("D20").Value = "26.585.98"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").Value = "'14.28"
$ws.Range("E21").Value = "  -3.15%  "
$ws.Range("D22").Value = "'5.022"
$ws.Range("E22").Value = "  -1.13%  "
$ws.Range("E23").Value = "  -0.32%  "
$ws.Range("E24").Value = "  -2.30%  "
$ws.Range("D25").Value = "'151.23"
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("D26").Value = "'17.95"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("E27").Value = "  +1.77%  "
$ws.Range("D28").Value = "'113.15"
$ws.Range("E28").Value = "  -2.44%  "
$ws.Range("D29").Value = "'4.862"
$ws.Range("E29").Value = "  -3.69%  "
$ws.Range("D30").Value = "'0.08697"
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("D31").Value = "'3.036"
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("E32").Value = "  -0.50%  "
$ws.Range("D33").Value = "'0.7317"
$ws.Range("E33").Value = "  -4.91%  "
$ws.Range("D34").Value = "'1.120"
$ws.Range("E34").Value = "  -4.15%  "
$ws.Range("D35").Value = "'2.662"
$ws.Range("E35").Value = "  -2.17%  "
$ws.Range("D36").Value = "'1.005"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("D37").Value = "'1.086"
$ws.Range("E37").Value = "  -1.85%  "
$ws.Range("D38").Value = "'0.01962"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'0.05132"
$ws.Range("E39").Value = "  -1.94%  "
$ws.Range("D40").Value = "'2.891"
$ws.Range("E40").Value = "  -1.87%  "
$ws.Range("D41").Value = "'6.987"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("D42").Value = "'0.5007"
$ws.Range("E42").Value = "  -2.21%  "
$ws.Range("E43").Value = "  -4.44%  "
$ws.Range("D44").Value = "'8.160"
$ws.Range("E44").Value = "  -4.42%  "
$ws.Range("D45").Value = "'1.007"
$ws.Range("E45").Value = "  +0.60%  "
$ws.Range("D46").Value = "'0.4605"
$ws.Range("E46").Value = "  -3.94%  "
$ws.Range("D47").Value = "'9.969"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("D48").Value = "'101.26"
$ws.Range("E48").Value = "  -1.48%  "
$ws.Range("D49").Value = "'1.589"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("D50").Value = "'0.06007"
$ws.Range("E50").Value = "  -3.09%  "
$ws.Range("D51").Value = "'64.56"
$ws.Range("E51").Value = "  -1.68%  "
